# Auto-generated script applying scheduled-runner market data updates
# to the Leve profit tables (columns H-N) across all job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2059.476
$ws.Range("I32").Value = 1557
$ws.Range("K32").Value = 1557
$ws.Range("M32").Value = -1231
$ws.Range("H70").Value = 2450
$ws.Range("J70").Value = 3000
$ws.Range("L70").Value = 9000
$ws.Range("N70").Value = -9540
$ws.Range("H73").Value = 2450
$ws.Range("J73").Value = 3000
$ws.Range("L73").Value = 9000
$ws.Range("N73").Value = -10872
$ws.Range("H132").Value = 22028.5
$ws.Range("I132").Value = 28927.027
$ws.Range("K132").Value = 86781.08099999999
$ws.Range("M132").Value = -84251.08099999999
$ws.Range("H138").Value = 6174.6875
$ws.Range("J138").Value = 3126.5757
$ws.Range("L138").Value = 9379.7271
$ws.Range("N138").Value = -19659.7271
$ws.Range("H141").Value = 2938
$ws.Range("I141").Value = 2999.5
$ws.Range("J141").Value = 2815
$ws.Range("K141").Value = 8998.5
$ws.Range("L141").Value = 8445
$ws.Range("M141").Value = -3818.5
$ws.Range("N141").Value = -18805

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 75218.5
$ws.Range("I45").Value = 103357
$ws.Range("J45").Value = 4872.25
$ws.Range("K45").Value = 103357
$ws.Range("L45").Value = 4872.25
$ws.Range("M45").Value = -102980
$ws.Range("N45").Value = -5626.25
$ws.Range("H74").Value = 379876.1
$ws.Range("I74").Value = 1760.129
$ws.Range("J74").Value = 798504.4399999999
$ws.Range("K74").Value = 1760.129
$ws.Range("L74").Value = 798504.4399999999
$ws.Range("M74").Value = -886.1289999999999
$ws.Range("N74").Value = -800252.4399999999
$ws.Range("H77").Value = 379876.1
$ws.Range("I77").Value = 1760.129
$ws.Range("J77").Value = 798504.4399999999
$ws.Range("K77").Value = 8800.645
$ws.Range("L77").Value = 3992522.2
$ws.Range("M77").Value = -4432.645
$ws.Range("N77").Value = -4001258.2
$ws.Range("H97").Value = 6315.8335
$ws.Range("I97").Value = 6315.8335
$ws.Range("K97").Value = 6315.8335
$ws.Range("M97").Value = -5819.8335
$ws.Range("H102").Value = 2737.7827
$ws.Range("I102").Value = 2725.8635
$ws.Range("K102").Value = 2725.8635
$ws.Range("M102").Value = -1103.8635
$ws.Range("H132").Value = 1808
$ws.Range("I132").Value = 1327.9584
$ws.Range("K132").Value = 3983.8752
$ws.Range("M132").Value = -1453.8752

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 2459.875
$ws.Range("I5").Value = 317.25
$ws.Range("J5").Value = 4602.5
$ws.Range("K5").Value = 317.25
$ws.Range("L5").Value = 4602.5
$ws.Range("M5").Value = -204.25
$ws.Range("N5").Value = -4828.5
$ws.Range("H60").Value = 42890
$ws.Range("J60").Value = 42890
$ws.Range("L60").Value = 42890
$ws.Range("N60").Value = -44088
$ws.Range("H94").Value = 4087.9
$ws.Range("I94").Value = 3647.5
$ws.Range("K94").Value = 3647.5
$ws.Range("M94").Value = -3196.5
$ws.Range("H107").Value = 12875.695
$ws.Range("J107").Value = 3931
$ws.Range("L107").Value = 3931
$ws.Range("N107").Value = -7771

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1785.7
$ws.Range("I134").Value = 1104.95
$ws.Range("J134").Value = 3147.2
$ws.Range("K134").Value = 3314.85
$ws.Range("L134").Value = 9441.599999999999
$ws.Range("M134").Value = -779.8500000000004
$ws.Range("N134").Value = -14511.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 6556127
$ws.Range("I4").Value = 10636481
$ws.Range("J4").Value = 144142.28
$ws.Range("K4").Value = 31909443
$ws.Range("L4").Value = 432426.84
$ws.Range("M4").Value = -31909331
$ws.Range("N4").Value = -432650.84
$ws.Range("H34").Value = 2185.1177
$ws.Range("I34").Value = 238.14285
$ws.Range("K34").Value = 714.4285500000001
$ws.Range("M34").Value = -630.4285500000001
$ws.Range("H64").Value = 7308.1665
$ws.Range("J64").Value = 7770
$ws.Range("L64").Value = 23310
$ws.Range("N64").Value = -23850
$ws.Range("H67").Value = 7308.1665
$ws.Range("J67").Value = 7770
$ws.Range("L67").Value = 23310
$ws.Range("N67").Value = -25182
$ws.Range("H113").Value = 1163.4166
$ws.Range("I113").Value = 1296.6666
$ws.Range("J113").Value = 1119
$ws.Range("K113").Value = 3889.9998
$ws.Range("L113").Value = 3357
$ws.Range("M113").Value = -1719.9998
$ws.Range("N113").Value = -7697
$ws.Range("H116").Value = 4323.25
$ws.Range("I116").Value = 4323.25
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 12969.75
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -9527.75
$ws.Range("N116").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 12321.3
$ws.Range("I70").Value = 11339.125
$ws.Range("J70").Value = 16250
$ws.Range("K70").Value = 11339.125
$ws.Range("L70").Value = 16250
$ws.Range("M70").Value = -11069.125
$ws.Range("N70").Value = -16790
$ws.Range("H73").Value = 12321.3
$ws.Range("I73").Value = 11339.125
$ws.Range("J73").Value = 16250
$ws.Range("K73").Value = 11339.125
$ws.Range("L73").Value = 16250
$ws.Range("M73").Value = -10403.125
$ws.Range("N73").Value = -18122
$ws.Range("H80").Value = 145215.16
$ws.Range("J80").Value = 82162.5
$ws.Range("L80").Value = 82162.5
$ws.Range("N80").Value = -84158.5
$ws.Range("H83").Value = 145215.16
$ws.Range("J83").Value = 82162.5
$ws.Range("L83").Value = 410812.5
$ws.Range("N83").Value = -420796.5
$ws.Range("H97").Value = 834.8946999999999
$ws.Range("I97").Value = 793.38464
$ws.Range("J97").Value = 924.8333
$ws.Range("K97").Value = 793.38464
$ws.Range("L97").Value = 924.8333
$ws.Range("M97").Value = -297.38464
$ws.Range("N97").Value = -1916.8333
$ws.Range("H113").Value = 2000.0968
$ws.Range("I113").Value = 1977.6897
$ws.Range("J113").Value = 2325
$ws.Range("K113").Value = 1977.6897
$ws.Range("L113").Value = 2325
$ws.Range("M113").Value = 192.3103000000001
$ws.Range("N113").Value = -6665
$ws.Range("H122").Value = 4807.826
$ws.Range("I122").Value = 2955.611
$ws.Range("K122").Value = 8866.832999999999
$ws.Range("M122").Value = -6416.832999999999
$ws.Range("H131").Value = 49887.5
$ws.Range("J131").Value = 49887.5
$ws.Range("L131").Value = 49887.5
$ws.Range("N131").Value = -59967.5
$ws.Range("H132").Value = 15315733
$ws.Range("I132").Value = 526
$ws.Range("J132").Value = 25061774
$ws.Range("K132").Value = 1578
$ws.Range("L132").Value = 75185322
$ws.Range("M132").Value = 952
$ws.Range("N132").Value = -75190382

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 625
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 625
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 625
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -965
$ws.Range("H40").Value = 2529634.8
$ws.Range("I40").Value = 3476299
$ws.Range("K40").Value = 3476299
$ws.Range("M40").Value = -3476163
$ws.Range("H61").Value = 4880052.5
$ws.Range("I61").Value = 5884161.5
$ws.Range("K61").Value = 5884161.5
$ws.Range("M61").Value = -5883959.5
$ws.Range("H82").Value = 4245.3335
$ws.Range("I82").Value = 1320.4445
$ws.Range("J82").Value = 8632.666999999999
$ws.Range("K82").Value = 1320.4445
$ws.Range("L82").Value = 8632.666999999999
$ws.Range("M82").Value = -959.4445000000001
$ws.Range("N82").Value = -9354.666999999999
$ws.Range("H85").Value = 4245.3335
$ws.Range("I85").Value = 1320.4445
$ws.Range("J85").Value = 8632.666999999999
$ws.Range("K85").Value = 1320.4445
$ws.Range("L85").Value = 8632.666999999999
$ws.Range("M85").Value = -72.44450000000006
$ws.Range("N85").Value = -11128.667
$ws.Range("H93").Value = 2892.3845
$ws.Range("I93").Value = 1475.125
$ws.Range("K93").Value = 1475.125
$ws.Range("M93").Value = -227.125
$ws.Range("H113").Value = 4880052.5
$ws.Range("I113").Value = 5884161.5
$ws.Range("K113").Value = 5884161.5
$ws.Range("M113").Value = -5881991.5
$ws.Range("H122").Value = 4254.75
$ws.Range("I122").Value = 2621.8333
$ws.Range("J122").Value = 5887.6665
$ws.Range("K122").Value = 7865.499899999999
$ws.Range("L122").Value = 17662.9995
$ws.Range("M122").Value = -5415.499899999999
$ws.Range("N122").Value = -22562.9995
$ws.Range("H132").Value = 5933.4614
$ws.Range("I132").Value = 2275.037
$ws.Range("J132").Value = 14164.917
$ws.Range("K132").Value = 6825.110999999999
$ws.Range("L132").Value = 42494.751
$ws.Range("M132").Value = -4295.110999999999
$ws.Range("N132").Value = -47554.751
$ws.Range("H136").Value = 4704.294
$ws.Range("I136").Value = 4132.778
$ws.Range("J136").Value = 5347.25
$ws.Range("K136").Value = 12398.334
$ws.Range("L136").Value = 16041.75
$ws.Range("M136").Value = -9848.334000000001
$ws.Range("N136").Value = -21141.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 496
$ws.Range("I113").Value = 159.33333
$ws.Range("K113").Value = 477.99999
$ws.Range("M113").Value = 1692.00001
$ws.Range("H122").Value = 2061.647
$ws.Range("I122").Value = 2096.75
$ws.Range("K122").Value = 6290.25
$ws.Range("M122").Value = -3840.25
$ws.Range("H126").Value = 2101.7
$ws.Range("I126").Value = 1981.4073
$ws.Range("J126").Value = 3184.3333
$ws.Range("K126").Value = 5944.2219
$ws.Range("L126").Value = 9552.999899999999
$ws.Range("M126").Value = -3474.2219
$ws.Range("N126").Value = -14492.9999
$ws.Range("H132").Value = 1849
$ws.Range("I132").Value = 1622.4546
$ws.Range("J132").Value = 2347.4
$ws.Range("K132").Value = 4867.3638
$ws.Range("L132").Value = 7042.200000000001
$ws.Range("M132").Value = -2337.3638
$ws.Range("N132").Value = -12102.2
$ws.Range("H136").Value = 21811.898
$ws.Range("I136").Value = 29492.572
$ws.Range("K136").Value = 88477.716
$ws.Range("M136").Value = -85927.716
